# Apply the "想去人数" (F) and "最低票价" (G) updates described by the diff.
# Sheet 1 = 展览 (Exhibition), Sheet 2 = 演出 (Performance),
# Sheet 3 = 本地生活 (Local life, unchanged), Sheet 4 = 全部类型 (All types, combined view).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (index 1) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G2").Value = 60
$ws1.Range("F4").Value = 271
$ws1.Range("F6").Value = 10129
$ws1.Range("F8").Value = 923
$ws1.Range("F9").Value = 1260
$ws1.Range("F10").Value = 6090
$ws1.Range("F15").Value = 3115
$ws1.Range("F18").Value = 608
$ws1.Range("F19").Value = 116
$ws1.Range("F20").Value = 25
$ws1.Range("F22").Value = 30
$ws1.Range("F23").Value = 1555

# ---- Sheet "演出" (index 2) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = "不可售"

# ---- Sheet "全部类型" (index 4) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G2").Value = 60
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F5").Value = 271
$ws4.Range("F7").Value = 10129
$ws4.Range("F9").Value = 923
$ws4.Range("F10").Value = 1260
$ws4.Range("F11").Value = 6090
$ws4.Range("F16").Value = 3115
$ws4.Range("F19").Value = 608
$ws4.Range("F20").Value = 116
$ws4.Range("F21").Value = 25
$ws4.Range("F23").Value = 30
$ws4.Range("F24").Value = 1555
